$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 41
$ws_ALC.Range("H41").Value = 1437.875
$ws_ALC.Range("I41").Value = 600.75
$ws_ALC.Range("J41").Value = 5623.5
$ws_ALC.Range("K41").Value = 600.75
$ws_ALC.Range("L41").Value = 5623.5
$ws_ALC.Range("M41").Value = -160.75
$ws_ALC.Range("N41").Value = -6503.5

# ALC row 64
$ws_ALC.Range("H64").Value = 7848.769
$ws_ALC.Range("J64").Value = 8230.817999999999
$ws_ALC.Range("L64").Value = 8230.817999999999
$ws_ALC.Range("N64").Value = -8726.817999999999

# ALC row 67
$ws_ALC.Range("H67").Value = 7848.769
$ws_ALC.Range("J67").Value = 8230.817999999999
$ws_ALC.Range("L67").Value = 8230.817999999999
$ws_ALC.Range("N67").Value = -9946.817999999999

# ARM row 32
$ws_ARM.Range("H32").Value = 4052.5898
$ws_ARM.Range("I32").Value = 3209.111
$ws_ARM.Range("K32").Value = 3209.111
$ws_ARM.Range("M32").Value = -2922.111

# ARM row 43
$ws_ARM.Range("H43").Value = 35035.4
$ws_ARM.Range("I43").Value = 69888
$ws_ARM.Range("K43").Value = 69888
$ws_ARM.Range("M43").Value = -69575

# ARM row 74
$ws_ARM.Range("H74").Value = 37041060
$ws_ARM.Range("I74").Value = 47621508
$ws_ARM.Range("K74").Value = 47621508
$ws_ARM.Range("M74").Value = -47620634

# ARM row 77
$ws_ARM.Range("H77").Value = 37041060
$ws_ARM.Range("I77").Value = 47621508
$ws_ARM.Range("K77").Value = 238107540
$ws_ARM.Range("M77").Value = -238103172

# BSM row 81
$ws_BSM.Range("H81").Value = 69778.8
$ws_BSM.Range("J81").Value = 69778.8
$ws_BSM.Range("L81").Value = 69778.8
$ws_BSM.Range("N81").Value = -71900.8

# BSM row 84
$ws_BSM.Range("H84").Value = 69778.8
$ws_BSM.Range("J84").Value = 69778.8
$ws_BSM.Range("L84").Value = 209336.4
$ws_BSM.Range("N84").Value = -219944.4

# BSM row 134
$ws_BSM.Range("H134").Value = 4734.2666
$ws_BSM.Range("I134").Value = 2273.4546
$ws_BSM.Range("J134").Value = 11501.5
$ws_BSM.Range("K134").Value = 6820.3638
$ws_BSM.Range("L134").Value = 34504.5
$ws_BSM.Range("M134").Value = -4285.3638
$ws_BSM.Range("N134").Value = -39574.5

# BSM row 139
$ws_BSM.Range("H139").Value = 46663.332
$ws_BSM.Range("J139").Value = 46663.332
$ws_BSM.Range("L139").Value = 46663.332
$ws_BSM.Range("N139").Value = -56943.332

# CRP row 22
$ws_CRP.Range("H22").Value = 2014.7693
$ws_CRP.Range("I22").Value = 1622.1111
$ws_CRP.Range("J22").Value = 2898.25
$ws_CRP.Range("K22").Value = 1622.1111
$ws_CRP.Range("L22").Value = 2898.25
$ws_CRP.Range("M22").Value = -1272.1111
$ws_CRP.Range("N22").Value = -3598.25

# CRP row 31
$ws_CRP.Range("H31").Value = 21738.178
$ws_CRP.Range("I31").Value = 3473.7317
$ws_CRP.Range("J31").Value = 57397.332
$ws_CRP.Range("K31").Value = 3473.7317
$ws_CRP.Range("L31").Value = 57397.332
$ws_CRP.Range("M31").Value = -3178.7317
$ws_CRP.Range("N31").Value = -57987.332

# CRP row 34
$ws_CRP.Range("H34").Value = 21738.178
$ws_CRP.Range("I34").Value = 3473.7317
$ws_CRP.Range("J34").Value = 57397.332
$ws_CRP.Range("K34").Value = 3473.7317
$ws_CRP.Range("L34").Value = 57397.332
$ws_CRP.Range("M34").Value = -3271.7317
$ws_CRP.Range("N34").Value = -57801.332

# CRP row 58
$ws_CRP.Range("H58").Value = 5801.8086
$ws_CRP.Range("I58").Value = 5003.3237
$ws_CRP.Range("J58").Value = 7890.154
$ws_CRP.Range("K58").Value = 5003.3237
$ws_CRP.Range("L58").Value = 7890.154
$ws_CRP.Range("M58").Value = -4800.3237
$ws_CRP.Range("N58").Value = -8296.154

# CRP row 62
$ws_CRP.Range("H62").Value = 7426.273
$ws_CRP.Range("I62").Value = 3240.2856
$ws_CRP.Range("J62").Value = 14751.75
$ws_CRP.Range("K62").Value = 3240.2856
$ws_CRP.Range("L62").Value = 14751.75
$ws_CRP.Range("M62").Value = -2616.2856
$ws_CRP.Range("N62").Value = -15999.75

# CRP row 65
$ws_CRP.Range("H65").Value = 7426.273
$ws_CRP.Range("I65").Value = 3240.2856
$ws_CRP.Range("J65").Value = 14751.75
$ws_CRP.Range("K65").Value = 16201.428
$ws_CRP.Range("L65").Value = 73758.75
$ws_CRP.Range("M65").Value = -13081.428
$ws_CRP.Range("N65").Value = -79998.75

# CRP row 99
$ws_CRP.Range("H99").Value = 0
$ws_CRP.Range("I99").Value = 0
$ws_CRP.Range("K99").Value = 0
$ws_CRP.Range("M99").ClearContents() | Out-Null

# CRP row 107
$ws_CRP.Range("H107").Value = 1008.7407
$ws_CRP.Range("I107").Value = 815.6667
$ws_CRP.Range("K107").Value = 815.6667
$ws_CRP.Range("M107").Value = 1104.3333

# CRP row 126
$ws_CRP.Range("H126").Value = 0
$ws_CRP.Range("I126").Value = 0
$ws_CRP.Range("K126").Value = 0
$ws_CRP.Range("M126").ClearContents() | Out-Null

# CRP row 132
$ws_CRP.Range("H132").Value = 4785.824
$ws_CRP.Range("I132").Value = 4631.6313
$ws_CRP.Range("J132").Value = 5302.8237
$ws_CRP.Range("K132").Value = 13894.8939
$ws_CRP.Range("L132").Value = 15908.4711
$ws_CRP.Range("M132").Value = -11364.8939
$ws_CRP.Range("N132").Value = -20968.4711

# CRP row 134
$ws_CRP.Range("H134").Value = 2597.1177
$ws_CRP.Range("I134").Value = 1609.1333
$ws_CRP.Range("K134").Value = 4827.3999
$ws_CRP.Range("M134").Value = -2292.3999

# CRP row 136
$ws_CRP.Range("H136").Value = 5801.8086
$ws_CRP.Range("I136").Value = 5003.3237
$ws_CRP.Range("J136").Value = 7890.154
$ws_CRP.Range("K136").Value = 15009.9711
$ws_CRP.Range("L136").Value = 23670.462
$ws_CRP.Range("M136").Value = -12459.9711
$ws_CRP.Range("N136").Value = -28770.462

# CRP row 138
$ws_CRP.Range("H138").Value = 56664.332
$ws_CRP.Range("J138").Value = 56664.332
$ws_CRP.Range("L138").Value = 56664.332
$ws_CRP.Range("N138").Value = -66944.33199999999

# CUL row 12
$ws_CUL.Range("H12").Value = 94.36364
$ws_CUL.Range("J12").Value = 73.625
$ws_CUL.Range("L12").Value = 220.875
$ws_CUL.Range("N12").Value = -566.875

# CUL row 15
$ws_CUL.Range("H15").Value = 380
$ws_CUL.Range("I15").Value = 55.6
$ws_CUL.Range("K15").Value = 166.8
$ws_CUL.Range("M15").Value = -26.80000000000001

# CUL row 39
$ws_CUL.Range("H39").Value = 3033.3333
$ws_CUL.Range("J39").Value = 2550
$ws_CUL.Range("L39").Value = 7650
$ws_CUL.Range("N39").Value = -8238

# CUL row 132
$ws_CUL.Range("H132").Value = 4092.923
$ws_CUL.Range("I132").Value = 2634
$ws_CUL.Range("J132").Value = 7375.5
$ws_CUL.Range("K132").Value = 23706
$ws_CUL.Range("L132").Value = 66379.5
$ws_CUL.Range("M132").Value = -21176
$ws_CUL.Range("N132").Value = -71439.5

# GSM row 126
$ws_GSM.Range("H126").Value = 6317.778
$ws_GSM.Range("I126").Value = 1798.6666
$ws_GSM.Range("J126").Value = 8577.333000000001
$ws_GSM.Range("K126").Value = 5395.9998
$ws_GSM.Range("L126").Value = 25731.999
$ws_GSM.Range("M126").Value = -2925.9998
$ws_GSM.Range("N126").Value = -30671.999

# GSM row 131
$ws_GSM.Range("H131").Value = 55995
$ws_GSM.Range("J131").Value = 55995
$ws_GSM.Range("L131").Value = 55995
$ws_GSM.Range("N131").Value = -66075

# GSM row 140
$ws_GSM.Range("H140").Value = 73113.39999999999
$ws_GSM.Range("J140").Value = 73113.39999999999
$ws_GSM.Range("L140").Value = 73113.39999999999
$ws_GSM.Range("N140").Value = -83473.39999999999

# LTW row 7
$ws_LTW.Range("H7").Value = 10801.25
$ws_LTW.Range("I7").Value = 2912.7144
$ws_LTW.Range("J7").Value = 21845.2
$ws_LTW.Range("K7").Value = 2912.7144
$ws_LTW.Range("L7").Value = 21845.2
$ws_LTW.Range("M7").Value = -2800.7144
$ws_LTW.Range("N7").Value = -22069.2

# LTW row 40
$ws_LTW.Range("H40").Value = 13757
$ws_LTW.Range("I40").Value = 12882.333
$ws_LTW.Range("K40").Value = 12882.333
$ws_LTW.Range("M40").Value = -12746.333

# LTW row 46
$ws_LTW.Range("H46").Value = 6204.5454
$ws_LTW.Range("J46").Value = 6761.1113
$ws_LTW.Range("L46").Value = 6761.1113
$ws_LTW.Range("N46").Value = -7137.1113

# LTW row 122
$ws_LTW.Range("H122").Value = 9684.333000000001
$ws_LTW.Range("I122").Value = 8742.333000000001
$ws_LTW.Range("J122").Value = 11568.333
$ws_LTW.Range("K122").Value = 26226.999
$ws_LTW.Range("L122").Value = 34704.999
$ws_LTW.Range("M122").Value = -23776.999
$ws_LTW.Range("N122").Value = -39604.999

# LTW row 126
$ws_LTW.Range("H126").Value = 10801.25
$ws_LTW.Range("I126").Value = 2912.7144
$ws_LTW.Range("J126").Value = 21845.2
$ws_LTW.Range("K126").Value = 8738.143199999999
$ws_LTW.Range("L126").Value = 65535.60000000001
$ws_LTW.Range("M126").Value = -6268.143199999999
$ws_LTW.Range("N126").Value = -70475.60000000001

# LTW row 132
$ws_LTW.Range("H132").Value = 3108.9656
$ws_LTW.Range("J132").Value = 6803.75
$ws_LTW.Range("L132").Value = 20411.25
$ws_LTW.Range("N132").Value = -25471.25

# LTW row 133
$ws_LTW.Range("H133").Value = 59326
$ws_LTW.Range("J133").Value = 59326
$ws_LTW.Range("L133").Value = 59326
$ws_LTW.Range("N133").Value = -64386

# LTW row 136
$ws_LTW.Range("H136").Value = 3675.36
$ws_LTW.Range("I136").Value = 1414.5264
$ws_LTW.Range("K136").Value = 4243.5792
$ws_LTW.Range("M136").Value = -1693.5792

# WVR row 122
$ws_WVR.Range("H122").Value = 5431.7393
$ws_WVR.Range("I122").Value = 1506.7
$ws_WVR.Range("J122").Value = 12791.1875
$ws_WVR.Range("K122").Value = 4520.1
$ws_WVR.Range("L122").Value = 38373.5625
$ws_WVR.Range("M122").Value = -2070.1
$ws_WVR.Range("N122").Value = -43273.5625

# WVR row 135
$ws_WVR.Range("H135").Value = 64262.5
$ws_WVR.Range("J135").Value = 64262.5
$ws_WVR.Range("L135").Value = 64262.5
$ws_WVR.Range("N135").Value = -74402.5

# WVR row 136
$ws_WVR.Range("H136").Value = 4519.4
$ws_WVR.Range("I136").Value = 3318.638
$ws_WVR.Range("K136").Value = 9955.914000000001
$ws_WVR.Range("M136").Value = -7405.914000000001
